$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix lat/lon issue: the D column (Lat_D, degrees) should be negative
# (southern hemisphere), and the G column (computed Lat) formula needs to
# subtract minutes/seconds instead of adding them so the sign is respected.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 4).Value = -13
    $ws.Cells.Item($r, 7).Formula = "=D$r-(E$r/60)-(F$r/3600)"
}

# Select the full used range, matching the saved selection state.
$ws.Range("A1:N9").Select() | Out-Null
